$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 8.191447666666667
$ws.Range("H2").Value = 24.574343
$ws.Range("I2").Value = 0.185794284429433
$ws.Range("J2").Value = 0.185794284429433
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.383699
$ws.Range("N2").Value = 4.151097
$ws.Range("O2").Value = 0.08080976933214185
$ws.Range("P2").Value = 0.08080976933214185
$ws.Range("Q2").Value = 11.334497944919
$ws.Range("R2").Value = 102.010481504271
$ws.Range("S2").Value = 0.01501399326797283
$ws.Range("T2").Value = 0.01501399326797283

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 8.191447666666667
$ws.Range("H3").Value = 24.574343
$ws.Range("I3").Value = 0.185794284429433
$ws.Range("J3").Value = 0.185794284429433
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 2.462094
$ws.Range("N3").Value = 7.386282
$ws.Range("O3").Value = 0.1437893994387872
$ws.Range("P3").Value = 0.1437893994387872
$ws.Range("Q3").Value = 20.168114151414
$ws.Range("R3").Value = 181.513027362726
$ws.Range("S3").Value = 0.02671524857726738
$ws.Range("T3").Value = 0.02671524857726739

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 8.191447666666667
$ws.Range("H4").Value = 24.574343
$ws.Range("I4").Value = 0.185794284429433
$ws.Range("J4").Value = 0.185794284429433
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 4.336036333333333
$ws.Range("N4").Value = 13.008109
$ws.Range("O4").Value = 0.2532299986575496
$ws.Range("P4").Value = 0.2532299986575496
$ws.Range("Q4").Value = 35.51841470526522
$ws.Range("R4").Value = 319.665732347387
$ws.Range("S4").Value = 0.04704868639664571
$ws.Range("T4").Value = 0.04704868639664571

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 8.191447666666667
$ws.Range("H5").Value = 24.574343
$ws.Range("I5").Value = 0.185794284429433
$ws.Range("J5").Value = 0.185794284429433
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 8.941088000000001
$ws.Range("N5").Value = 26.823264
$ws.Range("O5").Value = 0.5221708325715213
$ws.Range("P5").Value = 0.5221708325715213
$ws.Range("Q5").Value = 73.24045443506134
$ws.Range("R5").Value = 659.164089915552
$ws.Range("S5").Value = 0.09701635618754705
$ws.Range("T5").Value = 0.09701635618754705

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 17.317702
$ws.Range("H6").Value = 51.95310600000001
$ws.Range("I6").Value = 0.3927913821808575
$ws.Range("J6").Value = 0.3927913821808576
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.383699
$ws.Range("N6").Value = 4.151097
$ws.Range("O6").Value = 0.08080976933214185
$ws.Range("P6").Value = 0.08080976933214185
$ws.Range("Q6").Value = 23.962486939698
$ws.Range("R6").Value = 215.662382457282
$ws.Range("S6").Value = 0.03174138098968827
$ws.Range("T6").Value = 0.03174138098968827

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 17.317702
$ws.Range("H7").Value = 51.95310600000001
$ws.Range("I7").Value = 0.3927913821808575
$ws.Range("J7").Value = 0.3927913821808576
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 2.462094
$ws.Range("N7").Value = 7.386282
$ws.Range("O7").Value = 0.1437893994387872
$ws.Range("P7").Value = 0.1437893994387872
$ws.Range("Q7").Value = 42.637810187988
$ws.Range("R7").Value = 383.7402916918921
$ws.Range("S7").Value = 0.05647923694851665
$ws.Range("T7").Value = 0.05647923694851667

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 17.317702
$ws.Range("H8").Value = 51.95310600000001
$ws.Range("I8").Value = 0.3927913821808575
$ws.Range("J8").Value = 0.3927913821808576
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 4.336036333333333
$ws.Range("N8").Value = 13.008109
$ws.Range("O8").Value = 0.2532299986575496
$ws.Range("P8").Value = 0.2532299986575496
$ws.Range("Q8").Value = 75.09018508183934
$ws.Range("R8").Value = 675.8116657365541
$ws.Range("S8").Value = 0.09946656118235561
$ws.Range("T8").Value = 0.09946656118235563

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 17.317702
$ws.Range("H9").Value = 51.95310600000001
$ws.Range("I9").Value = 0.3927913821808575
$ws.Range("J9").Value = 0.3927913821808576
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 8.941088000000001
$ws.Range("N9").Value = 26.823264
$ws.Range("O9").Value = 0.5221708325715213
$ws.Range("P9").Value = 0.5221708325715213
$ws.Range("Q9").Value = 154.839097539776
$ws.Range("R9").Value = 1393.551877857984
$ws.Range("S9").Value = 0.205104203060297
$ws.Range("T9").Value = 0.205104203060297

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 12.35128266666667
$ws.Range("H10").Value = 37.053848
$ws.Range("I10").Value = 0.2801455637905346
$ws.Range("J10").Value = 0.2801455637905346
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.383699
$ws.Range("N10").Value = 4.151097
$ws.Range("O10").Value = 0.08080976933214185
$ws.Range("P10").Value = 0.08080976933214185
$ws.Range("Q10").Value = 17.090457474584
$ws.Range("R10").Value = 153.814117271256
$ws.Range("S10").Value = 0.02263849838933593
$ws.Range("T10").Value = 0.02263849838933593

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 12.35128266666667
$ws.Range("H11").Value = 37.053848
$ws.Range("I11").Value = 0.2801455637905346
$ws.Range("J11").Value = 0.2801455637905346
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 2.462094
$ws.Range("N11").Value = 7.386282
$ws.Range("O11").Value = 0.1437893994387872
$ws.Range("P11").Value = 0.1437893994387872
$ws.Range("Q11").Value = 30.410018945904
$ws.Range("R11").Value = 273.690170513136
$ws.Range("S11").Value = 0.04028196237288142
$ws.Range("T11").Value = 0.04028196237288142

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 12.35128266666667
$ws.Range("H12").Value = 37.053848
$ws.Range("I12").Value = 0.2801455637905346
$ws.Range("J12").Value = 0.2801455637905346
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 4.336036333333333
$ws.Range("N12").Value = 13.008109
$ws.Range("O12").Value = 0.2532299986575496
$ws.Range("P12").Value = 0.2532299986575496
$ws.Range("Q12").Value = 53.5556104059369
$ws.Range("R12").Value = 482.000493653432
$ws.Range("S12").Value = 0.07094126074259556
$ws.Range("T12").Value = 0.07094126074259556

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 12.35128266666667
$ws.Range("H13").Value = 37.053848
$ws.Range("I13").Value = 0.2801455637905346
$ws.Range("J13").Value = 0.2801455637905346
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 8.941088000000001
$ws.Range("N13").Value = 26.823264
$ws.Range("O13").Value = 0.5221708325715213
$ws.Range("P13").Value = 0.5221708325715213
$ws.Range("Q13").Value = 110.4339052355413
$ws.Range("R13").Value = 993.9051471198721
$ws.Range("S13").Value = 0.1462838422857217
$ws.Range("T13").Value = 0.1462838422857217

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 6.228371
$ws.Range("H14").Value = 18.685113
$ws.Range("I14").Value = 0.1412687695991749
$ws.Range("J14").Value = 0.1412687695991749
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 1.383699
$ws.Range("N14").Value = 4.151097
$ws.Range("O14").Value = 0.08080976933214185
$ws.Range("P14").Value = 0.08080976933214185
$ws.Range("Q14").Value = 8.618190724329001
$ws.Range("R14").Value = 77.563716518961
$ws.Range("S14").Value = 0.01141589668514481
$ws.Range("T14").Value = 0.01141589668514481

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 6.228371
$ws.Range("H15").Value = 18.685113
$ws.Range("I15").Value = 0.1412687695991749
$ws.Range("J15").Value = 0.1412687695991749
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 2.462094
$ws.Range("N15").Value = 7.386282
$ws.Range("O15").Value = 0.1437893994387872
$ws.Range("P15").Value = 0.1437893994387872
$ws.Range("Q15").Value = 15.334834868874
$ws.Range("R15").Value = 138.013513819866
$ws.Range("S15").Value = 0.02031295154012176
$ws.Range("T15").Value = 0.02031295154012176

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 6.228371
$ws.Range("H16").Value = 18.685113
$ws.Range("I16").Value = 0.1412687695991749
$ws.Range("J16").Value = 0.1412687695991749
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 4.336036333333333
$ws.Range("N16").Value = 13.008109
$ws.Range("O16").Value = 0.2532299986575496
$ws.Range("P16").Value = 0.2532299986575496
$ws.Range("Q16").Value = 27.00644295347967
$ws.Range("R16").Value = 243.057986581317
$ws.Range("S16").Value = 0.03577349033595274
$ws.Range("T16").Value = 0.03577349033595274

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 6.228371
$ws.Range("H17").Value = 18.685113
$ws.Range("I17").Value = 0.1412687695991749
$ws.Range("J17").Value = 0.1412687695991749
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 8.941088000000001
$ws.Range("N17").Value = 26.823264
$ws.Range("O17").Value = 0.5221708325715213
$ws.Range("P17").Value = 0.5221708325715213
$ws.Range("Q17").Value = 55.688413207648
$ws.Range("R17").Value = 501.195718868832
$ws.Range("S17").Value = 0.07376643103795556
$ws.Range("T17").Value = 0.07376643103795556
